$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 (new): T_cd_subcool = 1 K, source HMW3  -- set first so shared string ordering matches
$ws.Range("A21").Value = "T_cd_subcool"
$ws.Range("B21").Value = 1
$ws.Range("C21").Value = '"K"'
$ws.Range("E21").Value = "HMW3"

# Row 18: fix typo in F18 definition string
$ws.Range("F18").Value = "différence T entre sortie de l'evaporateur I et etat 1"

# Row 17: T_max -> 838.15 ; add note "Limite technologique de la turbine" in F17
$ws.Range("B17").Formula = "=838.15"
$ws.Range("F17").Value = "Limite technologique de la turbine"

# Row 22 (new): p_3 = 310*10^5 Pa, source HMW3, definition "Limite technologique de la turbine"
$ws.Range("A22").Value = "p_3"
$ws.Range("B22").Formula = "=310*10^5"
$ws.Range("C22").Value = '"Pa"'
$ws.Range("E22").Value = "HMW3"
$ws.Range("F22").Value = "Limite technologique de la turbine"

# Row 14: T_cold_fluid_in -> 273.15+8
$ws.Range("B14").Formula = "=273.15+8"

# Row 15: T_cold_fluid_out -> 273.15+28
$ws.Range("B15").Formula = "=273.15+28"

# Row 20: p_hot_fluid ; add unit "Pa" in C20
$ws.Range("C20").Value = '"Pa"'

# Update the active selection to match the final edit location
$ws.Range("C15").Select() | Out-Null
